$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Avg" header label in E2 (new shared string, same formatting as the
# other header cells B2/C2).
$ws.Range("E2").Value = "Avg"
$ws.Range("C2").Copy()
$ws.Range("E2").PasteSpecial(-4122)   # xlPasteFormats

# New average formula in E3 (=AVERAGE(B3:B5)), matching the formatting of
# the adjacent data cell C3.
$ws.Range("E3").Formula = "=AVERAGE(B3:B5)"
$ws.Range("C3").Copy()
$ws.Range("E3").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0
